$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 66 ("「前置詞」حروف الجر ..." entry) entirely.
# This shifts all subsequent rows up by one (e.g. old row 67 becomes new row 66),
# and shrinks the used range from A1:C257 to A1:C256.
$ws.Rows.Item(66).Delete()
